$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.186.19'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').Value = '2.795.63'
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.89'
$ws.Range('E5').Value = '  +4.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '116.53'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.554'
$ws.Range('E7').Value = '  +4.14%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.17'
$ws.Range('E10').Value = '  +3.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0856'
$ws.Range('E11').Value = '  +3.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.14'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.132'
$ws.Range('E13').Value = '  +2.07%  '
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('D15').Value = '3.235.44'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '2.793.85'
$ws.Range('E16').Value = '  +1.82%  '
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '52.133.21'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('E19').Value = '  +6.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.11'
$ws.Range('E20').Value = '  +3.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.40'
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.23'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.26'
$ws.Range('E24').Value = '  -4.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.76'
$ws.Range('E25').Value = '  +6.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.65'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.25'
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('E29').Value = '  +0.28%  '
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.98'
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.25'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.73'
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0415'
$ws.Range('E34').Value = '  +17.11%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0823'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.98'
$ws.Range('E38').Value = '  -2.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.98'
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  +20.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.73'
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '128.02'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('E44').Value = '  +2.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.30'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.34'
$ws.Range('E46').Value = '  -2.60%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.075.84'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.37'
$ws.Range('E48').Value = '  +5.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.974'
$ws.Range('E49').Value = '  +17.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.53'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  -1.24%  '
